$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C106").NumberFormat = "@"

$ws.Range("C2").Value = "01/10/2012"
$ws.Range("D2").Value = 78.30711837063924
$ws.Range("C3").Value = "01/10/2014"
$ws.Range("D3").Value = 79.13402020982937
$ws.Range("C4").Value = "01/10/2016"
$ws.Range("D4").Value = 79.74215875027417
$ws.Range("C5").Value = "01/10/2018"
$ws.Range("D5").Value = 80.2262217702477
$ws.Range("C6").Value = "01/10/2020"
$ws.Range("D6").Value = 80.64086330391397
$ws.Range("C7").Value = "01/10/2022"
$ws.Range("D7").Value = 80.97820011179429
$ws.Range("C8").Value = "01/10/2024"
$ws.Range("D8").Value = 81.28709548394809
$ws.Range("C9").Value = "01/10/2012"
$ws.Range("D9").Value = 45.67100221818915
$ws.Range("C10").Value = "01/10/2014"
$ws.Range("D10").Value = 46.07028342328145
$ws.Range("C11").Value = "01/10/2016"
$ws.Range("D11").Value = 43.9519411205615
$ws.Range("C12").Value = "01/10/2018"
$ws.Range("D12").Value = 44.86731776579619
$ws.Range("C13").Value = "01/10/2020"
$ws.Range("D13").Value = 41.21151701850207
$ws.Range("C14").Value = "01/10/2022"
$ws.Range("D14").Value = 46.28749767095211
$ws.Range("C15").Value = "01/10/2024"
$ws.Range("D15").Value = 47.72956099175681
$ws.Range("C16").Value = "01/10/2012"
$ws.Range("D16").Value = 3.392821133293003
$ws.Range("C17").Value = "01/10/2014"
$ws.Range("D17").Value = 3.248539272386672
$ws.Range("C18").Value = "01/10/2016"
$ws.Range("D18").Value = 6.080959227938489
$ws.Range("C19").Value = "01/10/2018"
$ws.Range("D19").Value = 5.95439108160485
$ws.Range("C20").Value = "01/10/2020"
$ws.Range("D20").Value = 6.809290722506757
$ws.Range("C21").Value = "01/10/2022"
$ws.Range("D21").Value = 3.992919694428918
$ws.Range("C22").Value = "01/10/2024"
$ws.Range("D22").Value = 3.136824005921485
$ws.Range("C23").Value = "01/10/2012"
$ws.Range("D23").Value = 29.24329501915709
$ws.Range("C24").Value = "01/10/2014"
$ws.Range("D24").Value = 29.81519751416125
$ws.Range("C25").Value = "01/10/2016"
$ws.Range("D25").Value = 29.70925840177418
$ws.Range("C26").Value = "01/10/2018"
$ws.Range("D26").Value = 29.40451292284667
$ws.Range("C27").Value = "01/10/2020"
$ws.Range("D27").Value = 32.62052803658836
$ws.Range("C28").Value = "01/10/2022"
$ws.Range("D28").Value = 30.69824855599031
$ws.Range("C29").Value = "01/10/2024"
$ws.Range("D29").Value = 30.42117022890586
$ws.Range("C30").Value = "01/10/2012"
$ws.Range("D30").Value = 49.06331921758419
$ws.Range("C31").Value = "01/10/2014"
$ws.Range("D31").Value = 49.31832711378065
$ws.Range("C32").Value = "01/10/2016"
$ws.Range("D32").Value = 50.03290034849999
$ws.Range("C33").Value = "01/10/2018"
$ws.Range("D33").Value = 50.82170884740103
$ws.Range("C34").Value = "01/10/2020"
$ws.Range("D34").Value = 48.02080774100882
$ws.Range("C35").Value = "01/10/2022"
$ws.Range("D35").Value = 50.28041736538103
$ws.Range("C36").Value = "01/10/2024"
$ws.Range("D36").Value = 50.86592525504224
$ws.Range("C37").Value = "01/10/2012"
$ws.Range("D37").Value = 75.9836633482843
$ws.Range("C38").Value = "01/10/2014"
$ws.Range("D38").Value = 76.89063513805782
$ws.Range("C39").Value = "01/10/2016"
$ws.Range("D39").Value = 78.03952315615032
$ws.Range("C40").Value = "01/10/2018"
$ws.Range("D40").Value = 78.5403338091965
$ws.Range("C41").Value = "01/10/2020"
$ws.Range("D41").Value = 0
$ws.Range("C42").Value = "01/10/2022"
$ws.Range("D42").Value = 79.82142548745315
$ws.Range("C43").Value = "01/10/2024"
$ws.Range("D43").Value = 80.22009995207777
$ws.Range("C44").Value = "01/10/2012"
$ws.Range("D44").Value = 39.54162564270868
$ws.Range("C45").Value = "01/10/2014"
$ws.Range("D45").Value = 40.92170715882056
$ws.Range("C46").Value = "01/10/2016"
$ws.Range("D46").Value = 37.4596838860279
$ws.Range("C47").Value = "01/10/2018"
$ws.Range("D47").Value = 37.72361162516083
$ws.Range("C48").Value = "01/10/2020"
$ws.Range("D48").Value = 0
$ws.Range("C49").Value = "01/10/2022"
$ws.Range("D49").Value = 38.82527675595392
$ws.Range("C50").Value = "01/10/2024"
$ws.Range("D50").Value = 40.43609228452112
$ws.Range("C51").Value = "01/10/2012"
$ws.Range("D51").Value = 4.091456077015645
$ws.Range("C52").Value = "01/10/2014"
$ws.Range("D52").Value = 3.719991348857328
$ws.Range("C53").Value = "01/10/2016"
$ws.Range("D53").Value = 6.365045706445232
$ws.Range("C54").Value = "01/10/2018"
$ws.Range("D54").Value = 6.418865330724898
$ws.Range("C55").Value = "01/10/2020"
$ws.Range("D55").Value = 0
$ws.Range("C56").Value = "01/10/2022"
$ws.Range("D56").Value = 4.728597827400998
$ws.Range("C57").Value = "01/10/2024"
$ws.Range("D57").Value = 3.816663243650305
$ws.Range("C58").Value = "01/10/2012"
$ws.Range("D58").Value = 32.35058162855996
$ws.Range("C59").Value = "01/10/2014"
$ws.Range("D59").Value = 32.24893663037993
$ws.Range("C60").Value = "01/10/2016"
$ws.Range("D60").Value = 34.21657549136656
$ws.Range("C61").Value = "01/10/2018"
$ws.Range("D61").Value = 34.39961930946968
$ws.Range("C62").Value = "01/10/2020"
$ws.Range("D62").Value = 0
$ws.Range("C63").Value = "01/10/2022"
$ws.Range("D63").Value = 36.2658238778647
$ws.Range("C64").Value = "01/10/2024"
$ws.Range("D64").Value = 35.96905593208736
$ws.Range("C65").Value = "01/10/2012"
$ws.Range("D65").Value = 43.63308171972432
$ws.Range("C66").Value = "01/10/2014"
$ws.Range("D66").Value = 44.64169850767789
$ws.Range("C67").Value = "01/10/2016"
$ws.Range("D67").Value = 43.82294766478376
$ws.Range("C68").Value = "01/10/2018"
$ws.Range("D68").Value = 44.14071449972682
$ws.Range("C69").Value = "01/10/2020"
$ws.Range("D69").Value = 0
$ws.Range("C70").Value = "01/10/2022"
$ws.Range("D70").Value = 43.55560160958845
$ws.Range("C71").Value = "01/10/2024"
$ws.Range("D71").Value = 44.25104401999042
$ws.Range("C72").Value = "01/10/2012"
$ws.Range("D72").Value = 76.38888888888889
$ws.Range("C73").Value = "01/10/2014"
$ws.Range("D73").Value = 76.89514298683613
$ws.Range("C74").Value = "01/10/2016"
$ws.Range("D74").Value = 78.88641425389756
$ws.Range("C75").Value = "01/10/2018"
$ws.Range("D75").Value = 78.82764654418197
$ws.Range("C76").Value = "01/10/2020"
$ws.Range("D76").Value = 0
$ws.Range("C77").Value = "01/10/2022"
$ws.Range("D77").Value = 78.55329949238579
$ws.Range("C78").Value = "01/10/2024"
$ws.Range("D78").Value = 79.46688879633486
$ws.Range("C79").Value = "01/10/2012"
$ws.Range("D79").Value = 42.82407407407408
$ws.Range("C80").Value = "01/10/2014"
$ws.Range("D80").Value = 42.71448025419883
$ws.Range("C81").Value = "01/10/2016"
$ws.Range("D81").Value = 38.75278396436526
$ws.Range("C82").Value = "01/10/2018"
$ws.Range("D82").Value = 39.501312335958
$ws.Range("C83").Value = "01/10/2020"
$ws.Range("D83").Value = 0
$ws.Range("C84").Value = "01/10/2022"
$ws.Range("D84").Value = 40.82064297800338
$ws.Range("C85").Value = "01/10/2024"
$ws.Range("D85").Value = 42.44064972927946
$ws.Range("C86").Value = "01/10/2012"
$ws.Range("D86").Value = 4.583333333333333
$ws.Range("C87").Value = "01/10/2014"
$ws.Range("D87").Value = 4.221516114389469
$ws.Range("C88").Value = "01/10/2016"
$ws.Range("D88").Value = 6.948775055679287
$ws.Range("C89").Value = "01/10/2018"
$ws.Range("D89").Value = 7.042869641294838
$ws.Range("C90").Value = "01/10/2020"
$ws.Range("D90").Value = 0
$ws.Range("C91").Value = "01/10/2022"
$ws.Range("D91").Value = 5.541455160744501
$ws.Range("C92").Value = "01/10/2024"
$ws.Range("D92").Value = 3.915035401915869
$ws.Range("C93").Value = "01/10/2012"
$ws.Range("D93").Value = 28.98148148148148
$ws.Range("C94").Value = "01/10/2014"
$ws.Range("D94").Value = 29.95914661824785
$ws.Range("C95").Value = "01/10/2016"
$ws.Range("D95").Value = 33.2293986636971
$ws.Range("C96").Value = "01/10/2018"
$ws.Range("D96").Value = 32.23972003499563
$ws.Range("C97").Value = "01/10/2020"
$ws.Range("D97").Value = 0
$ws.Range("C98").Value = "01/10/2022"
$ws.Range("D98").Value = 32.1912013536379
$ws.Range("C99").Value = "01/10/2024"
$ws.Range("D99").Value = 33.11120366513953
$ws.Range("C100").Value = "01/10/2012"
$ws.Range("D100").Value = 47.40740740740741
$ws.Range("C101").Value = "01/10/2014"
$ws.Range("D101").Value = 46.93599636858829
$ws.Range("C102").Value = "01/10/2016"
$ws.Range("D102").Value = 45.70155902004454
$ws.Range("C103").Value = "01/10/2018"
$ws.Range("D103").Value = 46.58792650918635
$ws.Range("C104").Value = "01/10/2020"
$ws.Range("D104").Value = 0
$ws.Range("C105").Value = "01/10/2022"
$ws.Range("D105").Value = 46.36209813874789
$ws.Range("C106").Value = "01/10/2024"
$ws.Range("D106").Value = 46.39733444398168

$ws.Range("C2:C106").ClearFormats()
